$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 was blank (no styling at all); give it the same per-column
# formatting as the other data rows before filling in values. (Row 5 has
# no K column entry, unlike row 6, so only copy A:J formatting.)
$ws.Range("A6:J6").Copy()
$ws.Range("A5:J5").PasteSpecial(-4122) # xlPasteFormats

# Row 5 was blank; fill it in with the new "Crumpet" facility data (no
# rows are shifted - row 6 below already existed as a blank, styled row
# and simply gets populated in place).
$ws.Range("A5").Value = "Crumpet GEF"
$ws.Range("B5").Value = 20001371
$ws.Range("C5").Value = "Crumpet exporter"
$ws.Range("D5").Value = "GBP"
$ws.Range("E5").Value = 7000000
$ws.Range("F5").Value = 3938753.8
$ws.Range("G5").Value = 777
$ws.Range("H5").Value = 456
$ws.Range("I5").Value = "GBP"
$ws.Range("J5").Value = "GBP"

# Row 6 - Scone GEF / Scone exporter
$ws.Range("A6").Value = "Scone GEF"
$ws.Range("B6").Value = 20001371
$ws.Range("C6").Value = "Scone exporter"
$ws.Range("D6").Value = "GBP"
$ws.Range("E6").Value = 770000
$ws.Range("F6").Value = 761579.37
$ws.Range("G6").Value = 777
$ws.Range("H6").Value = 456.77
$ws.Range("I6").Value = "GBP"
$ws.Range("J6").Value = "GBP"

$ws.Range("D7").Select()
